# The workbook contains a roster of companies (one per row, rows 2-21).
# The first data row (row 2, the "xxx有限公司（或缩写）" placeholder/template
# row) was removed, causing every subsequent company's details (columns
# B through S) to shift up by one row, and the final row (row 21) to be
# dropped entirely. The sequential index numbers in column A (0, 1, 2, ...)
# stay put - they simply keep counting 0..18 for the remaining 19 companies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (row 2); everything below shifts up one row,
# and the sheet shrinks from 21 rows to 20 rows.
$ws.Rows(2).Delete()

# Restore the sequential index numbers in column A for the remaining data
# rows (A2:A20 = 0..18), since the row delete/shift also shifted column A,
# but the index column is meant to stay fixed in place.
for ($i = 0; $i -le 18; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}
